$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new entries in column A, extending the sheetData the same way
# rows 17-20 do (single cell, column-A default style).
$ws.Range("A21").Value = "SCRIPT/D79P11A/enter03.ssb"
$ws.Range("A22").Value = "SCRIPT/D38P12A/enter06.ssb"

# Match the updated active-cell selection recorded in the sheet view.
$ws.Range("C5").Select()
